$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "B3"   = 5.948
    "D6"   = -7.912000000000001
    "B14"  = 5.77
    "D18"  = -8.308000000000002
    "D19"  = -8.02
    "B21"  = 9.159000000000001
    "C22"  = -12.711
    "B23"  = 8.843999999999999
    "C24"  = -12.367
    "B25"  = 5.596
    "B26"  = 5.587000000000001
    "C28"  = -12.938
    "B29"  = 5.602
    "C36"  = -12.732
    "D44"  = -7.242
    "C45"  = -13.026
    "D47"  = -7.591000000000001
    "C48"  = -11.347
    "C49"  = -12.91
    "D51"  = -8.272000000000002
    "C52"  = -11.429
    "B53"  = 5.479000000000001
    "C53"  = -10.494
    "C54"  = -12.776
    "D55"  = -8.266999999999999
    "B57"  = 4.888
    "D57"  = -8.213000000000003
    "B59"  = 4.968999999999999
    "D64"  = -7.787999999999999
    "B69"  = 6.274
    "C70"  = -11.531
    "B79"  = 6.609
    "D80"  = -8.074999999999999
    "B83"  = 5.1
    "C86"  = -13.477
    "C87"  = -13.362
    "C89"  = -13.207
    "B91"  = 5.948
    "D92"  = -6.842000000000001
    "B93"  = 6.069
    "D94"  = -6.813
    "D96"  = -7.35
    "C101" = -12.143
    "D101" = -7.687
    "B103" = 5.592999999999999
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
